# Refresh the cryptos list (Price / Volume(1h) columns, plus the Filecoin/OKB
# row swap) to the latest scraped snapshot.
# Note: several "Price" strings look numeric (e.g. "0.430", "1.00", "25.60");
# a leading apostrophe forces Excel to keep them as literal text instead of
# silently parsing them into doubles and dropping the trailing zeros.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.936.16"
$ws.Range("E2").Value = "  -1.71%  "
$ws.Range("D3").Value = "2.992.63"
$ws.Range("E3").Value = "  -4.31%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "'562.18"
$ws.Range("E5").Value = "  -4.06%  "
$ws.Range("D6").Value = "'127.80"
$ws.Range("E6").Value = "  -5.08%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").Value = "2.991.31"
$ws.Range("E8").Value = "  -4.24%  "
$ws.Range("D9").Value = "'0.495"
$ws.Range("E9").Value = "  -2.10%  "
$ws.Range("D10").Value = "'0.134"
$ws.Range("E10").Value = "  -4.75%  "
$ws.Range("D11").Value = "'5.17"
$ws.Range("E11").Value = "  -0.91%  "
$ws.Range("D12").Value = "'0.427"
$ws.Range("E12").Value = "  -5.42%  "
$ws.Range("D13").Value = "'0.0000221"
$ws.Range("E13").Value = "  -4.48%  "
$ws.Range("D14").Value = "'32.58"
$ws.Range("E14").Value = "  -3.83%  "
$ws.Range("E15").Value = "  +0.33%  "
$ws.Range("D16").Value = "3.503.72"
$ws.Range("E16").Value = "  -3.69%  "
$ws.Range("D17").Value = "61.137.93"
$ws.Range("E17").Value = "  -1.40%  "
$ws.Range("D18").Value = "3.013.83"
$ws.Range("E18").Value = "  -3.58%  "
$ws.Range("D19").Value = "'6.18"
$ws.Range("E19").Value = "  -4.84%  "
$ws.Range("D20").Value = "'436.45"
$ws.Range("E20").Value = "  -3.00%  "
$ws.Range("D21").Value = "'13.05"
$ws.Range("E21").Value = "  -5.24%  "
$ws.Range("D22").Value = "'0.658"
$ws.Range("E22").Value = "  -5.84%  "
$ws.Range("D23").Value = "'7.11"
$ws.Range("E23").Value = "  -5.64%  "
$ws.Range("D24").Value = "'78.68"
$ws.Range("E24").Value = "  -5.08%  "
$ws.Range("D25").Value = "'12.42"
$ws.Range("E25").Value = "  -6.08%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("D27").Value = "'1.00"
$ws.Range("E27").Value = "  +0.10%  "
$ws.Range("D28").Value = "'2.48"
$ws.Range("E28").Value = "  -6.97%  "
$ws.Range("D29").Value = "'7.13"
$ws.Range("E29").Value = "  -6.46%  "
$ws.Range("D30").Value = "'6.15"
$ws.Range("E30").Value = "  -8.43%  "
$ws.Range("D31").Value = "'25.41"
$ws.Range("E31").Value = "  -6.01%  "
$ws.Range("D32").Value = "'1.86"
$ws.Range("E32").Value = "  -6.68%  "
$ws.Range("D33").Value = "'0.0930"
$ws.Range("E33").Value = "  -8.95%  "
$ws.Range("D34").Value = "'2.25"
$ws.Range("E34").Value = "  -3.97%  "
$ws.Range("D35").Value = "'0.948"
$ws.Range("E35").Value = "  -6.85%  "
$ws.Range("B36").Value = "OKB"
$ws.Range("C36").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D36").Value = "'49.83"
$ws.Range("E36").Value = "  -1.74%  "
$ws.Range("B37").Value = "Filecoin"
$ws.Range("C37").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D37").Value = "'5.50"
$ws.Range("E37").Value = "  -4.56%  "
$ws.Range("D38").Value = "0.0₃0671"
$ws.Range("E38").Value = "  -2.58%  "
$ws.Range("D39").Value = "'0.0360"
$ws.Range("E39").Value = "  -5.72%  "
$ws.Range("D40").Value = "'7.69"
$ws.Range("E40").Value = "  -3.86%  "
$ws.Range("D41").Value = "'0.107"
$ws.Range("E41").Value = "  -2.79%  "
$ws.Range("D42").Value = "'369.20"
$ws.Range("E42").Value = "  -5.89%  "
$ws.Range("D43").Value = "2.647.82"
$ws.Range("E43").Value = "  -3.67%  "
$ws.Range("D44").Value = "'2.41"
$ws.Range("E44").Value = "  -8.49%  "
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("D46").Value = "'0.234"
$ws.Range("E46").Value = "  -5.70%  "
$ws.Range("D47").Value = "'118.87"
$ws.Range("E47").Value = "  -4.41%  "
$ws.Range("D48").Value = "'1.95"
$ws.Range("E48").Value = "  -6.78%  "
$ws.Range("D49").Value = "'32.72"
$ws.Range("E49").Value = "  -4.06%  "
$ws.Range("E50").Value = "  -3.89%  "
$ws.Range("D51").Value = "'23.37"
$ws.Range("E51").Value = "  -6.71%  "
